# "add charge logic to service" -- bump the membership-fee / charge counts
# for a few checklist rows on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8  (1.T-Admin / 2) ) and Row 16 (2.E-Admin / 4) ) go from 1 charge to 3.
$ws.Range("D8").Value  = 3
$ws.Range("D16").Value = 3

# Rows 35-38 (6.Data consumer, items 7)-10)) switch from the numeric charge
# count "3" to the "1 done" status label used elsewhere in the sheet.
$ws.Range("D35").Value = "1 done"
$ws.Range("D36").Value = "1 done"
$ws.Range("D37").Value = "1 done"
$ws.Range("D38").Value = "1 done"

# Reflect the author's last on-screen selection (row scrolled down a bit,
# cursor resting on C26) without disturbing the existing frozen-pane split.
$ws.Range("C26").Select()
